$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "96.913.38"
$ws.Range("E2").Value = "  +0.13%  "

# Row 3
$ws.Range("D3").Value = "3.694.02"
$ws.Range("E3").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.67%  "

# Row 6
$ws.Range("E6").Value = "  +1.56%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "654.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.41%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.426"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.77%  "

# Row 9
$ws.Range("B9").Value = "USDC"
$ws.Range("C9").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.02%  "

# Row 10
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.80%  "

# Row 11
$ws.Range("D11").Value = "3.692.41"
$ws.Range("E11").Value = "  +0.18%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.79%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000299"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.92%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.40%  "

# Row 16
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "4.383.66"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "96.743.00"
$ws.Range("E17").Value = "  +0.35%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.72%  "

# Row 19
$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.693.54"
$ws.Range("E19").Value = "  +1.46%  "

# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.06%  "

# Row 21
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "

# Row 22
$ws.Range("B22").Value = "Stellar"
$ws.Range("C22").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.510"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.02%  "

# Row 23
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "522.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.43%  "

# Row 24
$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.70%  "

# Row 25
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000211"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.05%  "

# Row 26
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.44%  "

# Row 27
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.86%  "

# Row 28
$ws.Range("B28").Value = "Hedera"
$ws.Range("C28").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.196"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +18.18%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.47%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.10%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.97%  "

# Row 32
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.03%  "

# Row 33
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.187"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.38%  "

# Row 35
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.27%  "

# Row 36
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.39%  "

# Row 37
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "646.14"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.39%  "

# Row 38
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.599"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.04%  "

# Row 40
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.02%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.43%  "

# Row 42
$ws.Range("B42").Value = "ImmutableX"
$ws.Range("C42").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.89"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.66%  "

# Row 44
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.160"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.29%  "

# Row 45
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.954"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.452"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.53%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0459"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.63"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.01%  "

# Row 49
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.87%  "

# Row 50
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.80%  "

# Row 51
$ws.Range("B51").Value = "MantraDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.03%  "
